$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protein names in order (columns B..DQ on row 1, rows 2..121 in column A),
# stripping the Python tuple-repr wrapper, e.g. ("ANG_1",) -> ANG_1
$names = @(
    'ANG_1',
    'BDNF_1',
    'BLC_1',
    'BMP-4_1',
    'BMP-6_1',
    'CK b8-1_1',
    'CNTF_1',
    'EGF_1',
    'Eotaxin_1',
    'Eotaxin-2_1',
    'Eotaxin-3_1',
    'FGF-6_1',
    'FGF-7_1',
    'Fit-3 Ligand_1',
    'Fractalkine_1',
    'GCP-2_1',
    'GDNF_1',
    'GM-CSF_1',
    'I-309_1',
    'IFN-g_1',
    'IGF-1_1',
    'IGFBP-1_1',
    'IGFBP-2_1',
    'IGFBP-4_1',
    'IL-10_1',
    'IL-13_1',
    'IL-15_1',
    'IL-16_1',
    'IL-1a_1',
    'IL-1b_1',
    'IL-1ra_1',
    'IL-2_1',
    'IL-3_1',
    'IL-4_1',
    'IL-5_1',
    'IL-6_1',
    'IL-7_1',
    'LEPTIN(OB)_1',
    'LIGHT_1',
    'MCP-1_1',
    'MCP-2_1',
    'MCP-3_1',
    'MCP-4_1',
    'M-CSF_1',
    'MDC_1',
    'MIG_1',
    'MIP-1d_1',
    'MIP-3a_1',
    'NAP-2_1',
    'NT-3_1',
    'PARC_1',
    'PDGF-BB_1',
    'RANTES_1',
    'SCF_1',
    'SDF-1_1',
    'TARC_1',
    'TGF-b_1',
    'TGF-b3_1',
    'TNF-a_1',
    'TNF-b_1',
    'Acrp30_1',
    'AgRP(ART)_1',
    'ANG-2_1',
    'AR_1',
    'AXL_1',
    'bFGF',
    'b-NGF_1',
    'BTC_1',
    'CCL-28_1',
    'CTACK_1',
    'DTK_1',
    'EGF-R_1',
    'ENA-78_1',
    'FAS_1',
    'FGF-4_1',
    'FGF-9_1',
    'GCSF_1',
    'GITR_1',
    'GITR-Light_1',
    'GRO_1',
    'GRO-a_1',
    'HCC-4_1',
    'HGF_1',
    'ICAM-1_1',
    'ICAM-3_1',
    'IGF-1 SR',
    'IGFBP3_1',
    'IGFBP-6_1',
    'IL-1 RI_1',
    'IL-11_1',
    'IL-12 p40_1',
    'IL-12 p70_1',
    'IL-17_1',
    'IL-1R4 /ST2_1',
    'IL-2 Ra_1',
    'IL-6 R_1',
    'IL-8_1',
    'I-TAC_1',
    'Lymphotactin_1',
    'MIF_1',
    'MIP-1a_1',
    'MIP-1b_1',
    'MIP-3b_1',
    'MSP-a_1',
    'NT-4_1',
    'OSM_1',
    'OST_1',
    'PIGF_1',
    'spg130_1',
    'sTNF RI_1',
    'sTNF RII_1',
    'TECK_1',
    'TIMP-1_1',
    'TIMP-2_1',
    'TPO_1',
    'TRAIL R3_1',
    'TRAIL R4_1',
    'uPAR_1',
    'VEGF-B_1',
    'VEGF-D_1'
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $name = $names[$i]
    # Row 1 header cells start at column B (index 2)
    $ws.Cells.Item(1, $i + 2).Value = $name
    # Column A row headers start at row 2
    $ws.Cells.Item($i + 2, 1).Value = $name
}

# Zero out the redundant upper-triangle half of the symmetric matrix
# (row index < column index), keeping the lower-triangle values intact.
$zeroCells = @(
    'BB2',
    'BJ2',
    'E3',
    'DA3',
    'L4',
    'AH4',
    'AI4',
    'BI4',
    'F5',
    'I5',
    'H6',
    'AS6',
    'AV6',
    'H7',
    'AH7',
    'AT7',
    'CN7',
    'BA9',
    'U10',
    'M11',
    'P12',
    'T12',
    'CC12',
    'CE12',
    'AD13',
    'AH13',
    'BI13',
    'DA13',
    'O14',
    'AC14',
    'AJ14',
    'BK14',
    'BR14',
    'BV14',
    'P15',
    'T15',
    'AB15',
    'AE15',
    'AR15',
    'BF15',
    'S17',
    'AC17',
    'BE17',
    'W18',
    'Z18',
    'AN18',
    'BH18',
    'BU18',
    'DH18',
    'AA19',
    'AJ19',
    'V20',
    'AW20',
    'BQ20',
    'CV20',
    'AI21',
    'BI23',
    'AV24',
    'DM24',
    'AR25',
    'CZ26',
    'AG28',
    'CE29',
    'BU30',
    'CB30',
    'AV32',
    'AK33',
    'AM34',
    'BC34',
    'AO35',
    'AQ35',
    'AT35',
    'AW35',
    'AL36',
    'CP36',
    'AQ40',
    'AY41',
    'AQ42',
    'AR43',
    'AU43',
    'BI48',
    'BF49',
    'AZ50',
    'DJ50',
    'BF51',
    'DK52',
    'DK54',
    'BG55',
    'BF56',
    'BX57',
    'CM57',
    'DB57',
    'DC58',
    'BI59',
    'BN64',
    'BN65',
    'BP65',
    'CE65',
    'BO66',
    'BQ68',
    'CA68',
    'CO68',
    'BU69',
    'CF70',
    'CO70',
    'CG71',
    'CJ71',
    'CU71',
    'CH72',
    'DB72',
    'CC73',
    'CJ73',
    'CJ75',
    'CA77',
    'CC77',
    'CI77',
    'CO77',
    'CX77',
    'DJ78',
    'DK78',
    'DL81',
    'CE82',
    'CO82',
    'CI83',
    'CI86',
    'CV86',
    'CL87',
    'DL87',
    'CK88',
    'CX88',
    'DL88',
    'CM90',
    'CO90',
    'CQ90',
    'DQ90',
    'DC92',
    'DN92',
    'CZ93',
    'DB93',
    'DD93',
    'DE94',
    'DE96',
    'DF97',
    'DJ97',
    'DM97',
    'DL98',
    'CW100',
    'CY101',
    'DO101',
    'CY102',
    'DP103',
    'DQ104',
    'DO106',
    'DD107',
    'DG107',
    'DL107',
    'DE108',
    'DH110',
    'DI112',
    'DM112',
    'DN116',
    'DQ119'
)

foreach ($addr in $zeroCells) {
    $ws.Range($addr).Value = 0
}
